$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of columns A and B for rows 2-5 -------------------
# before: A2:A5 = 2,3,4,5  B2:B5 = 1,2,3,4
# after:  A2:A5 = 1,2,3,4  B2:B5 = 2,3,4,5
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 5

# --- Move the active selection from A4 to C6 ------------------------------
$ws.Range("C6").Select() | Out-Null

# --- Explicit portrait page setup ----------------------------------------
$ws.PageSetup.Orientation = 1
